# Fixed duplicate name bug
#
# The text scraper was creating two distinct "Found by"/"Allocated to"
# entries for the same person ("Josh" and "josh"). Normalise every
# lower-case "josh" entry to the existing canonical "Josh" shared string,
# and record the missing Solved Date for the "QS text scraper features
# duplicate names" row now that the bug has actually been fixed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the duplicate-name bug: "josh" -> "Josh" -------------------------
# (rows 6-9, columns C "Found by" and D "Allocated to")
$ws.Range("C6").Value = "Josh"
$ws.Range("D6").Value = "Josh"
$ws.Range("C7").Value = "Josh"
$ws.Range("D7").Value = "Josh"
$ws.Range("C8").Value = "Josh"
$ws.Range("D8").Value = "Josh"
$ws.Range("C9").Value = "Josh"
$ws.Range("D9").Value = "Josh"

# --- Record the Solved Date for row 5 --------------------------------------
# Copy the date formatting already used in column E (e.g. E6) onto E5 so the
# new value inherits the existing "Solved Date" number format instead of
# creating a brand new style, then write the date itself (2015-10-08).
$ws.Range("E6").Copy($ws.Range("E5"))
$ws.Range("E5").Value = Get-Date -Year 2015 -Month 10 -Day 8 -Hour 0 -Minute 0 -Second 0

# --- Update the sheet's remembered selection --------------------------------
$ws.Range("D10").Select()
